# Replace the 100 arithmetic expressions in the single 20x5 table,
# addressing each cell by (row, col) so the duplicate "98-4=" source
# values (row 7 col 1, row 14 col 1) resolve to their distinct targets.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(1,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "56-53="
$cellRange = $t.Cell(1,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "93-17="
$cellRange = $t.Cell(1,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "4+27="
$cellRange = $t.Cell(1,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "64-8="
$cellRange = $t.Cell(1,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "68-62="
$cellRange = $t.Cell(2,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "50+24="
$cellRange = $t.Cell(2,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "19-13="
$cellRange = $t.Cell(2,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "69+25="
$cellRange = $t.Cell(2,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "32+5="
$cellRange = $t.Cell(2,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "3+56="
$cellRange = $t.Cell(3,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "69-40="
$cellRange = $t.Cell(3,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "98-33="
$cellRange = $t.Cell(3,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "10+38="
$cellRange = $t.Cell(3,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "10+55="
$cellRange = $t.Cell(3,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "46+11="
$cellRange = $t.Cell(4,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "40+12="
$cellRange = $t.Cell(4,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "71-44="
$cellRange = $t.Cell(4,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "86-79="
$cellRange = $t.Cell(4,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "93-61="
$cellRange = $t.Cell(4,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "84-38="
$cellRange = $t.Cell(5,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "71-37="
$cellRange = $t.Cell(5,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "93-83="
$cellRange = $t.Cell(5,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "31+50="
$cellRange = $t.Cell(5,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "42+33="
$cellRange = $t.Cell(5,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "84-79="
$cellRange = $t.Cell(6,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "29+10="
$cellRange = $t.Cell(6,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "3+84="
$cellRange = $t.Cell(6,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "86-37="
$cellRange = $t.Cell(6,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "54+22="
$cellRange = $t.Cell(6,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "5+65="
$cellRange = $t.Cell(7,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "22-15="
$cellRange = $t.Cell(7,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "68-28="
$cellRange = $t.Cell(7,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "28-26="
$cellRange = $t.Cell(7,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "57-11="
$cellRange = $t.Cell(7,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "50+4="
$cellRange = $t.Cell(8,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "90-4="
$cellRange = $t.Cell(8,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "25+19="
$cellRange = $t.Cell(8,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "71-34="
$cellRange = $t.Cell(8,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "80+15="
$cellRange = $t.Cell(8,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "10+33="
$cellRange = $t.Cell(9,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "53+23="
$cellRange = $t.Cell(9,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "78-21="
$cellRange = $t.Cell(9,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "68-30="
$cellRange = $t.Cell(9,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "16+27="
$cellRange = $t.Cell(9,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "34-33="
$cellRange = $t.Cell(10,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "14+35="
$cellRange = $t.Cell(10,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "43-23="
$cellRange = $t.Cell(10,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "30+34="
$cellRange = $t.Cell(10,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "82-49="
$cellRange = $t.Cell(10,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "29-24="
$cellRange = $t.Cell(11,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "87-15="
$cellRange = $t.Cell(11,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "18+3="
$cellRange = $t.Cell(11,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "66+19="
$cellRange = $t.Cell(11,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "97-96="
$cellRange = $t.Cell(11,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "62-7="
$cellRange = $t.Cell(12,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "86-13="
$cellRange = $t.Cell(12,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "33+60="
$cellRange = $t.Cell(12,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "13+69="
$cellRange = $t.Cell(12,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "19+10="
$cellRange = $t.Cell(12,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "74+22="
$cellRange = $t.Cell(13,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "29+11="
$cellRange = $t.Cell(13,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "36-5="
$cellRange = $t.Cell(13,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "74+5="
$cellRange = $t.Cell(13,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "70+13="
$cellRange = $t.Cell(13,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "99-33="
$cellRange = $t.Cell(14,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "28+10="
$cellRange = $t.Cell(14,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "10+3="
$cellRange = $t.Cell(14,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "74-58="
$cellRange = $t.Cell(14,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "70-27="
$cellRange = $t.Cell(14,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "44+17="
$cellRange = $t.Cell(15,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "35+28="
$cellRange = $t.Cell(15,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "2-1="
$cellRange = $t.Cell(15,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "86-24="
$cellRange = $t.Cell(15,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "58-7="
$cellRange = $t.Cell(15,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "37+25="
$cellRange = $t.Cell(16,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "59+3="
$cellRange = $t.Cell(16,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "40+4="
$cellRange = $t.Cell(16,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "81-20="
$cellRange = $t.Cell(16,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "95-95="
$cellRange = $t.Cell(16,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "15+43="
$cellRange = $t.Cell(17,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "71-1="
$cellRange = $t.Cell(17,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "44+30="
$cellRange = $t.Cell(17,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "37+40="
$cellRange = $t.Cell(17,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "37+33="
$cellRange = $t.Cell(17,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "26+17="
$cellRange = $t.Cell(18,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "2+10="
$cellRange = $t.Cell(18,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "46-8="
$cellRange = $t.Cell(18,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "16+44="
$cellRange = $t.Cell(18,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "15+27="
$cellRange = $t.Cell(18,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "78-6="
$cellRange = $t.Cell(19,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "79+7="
$cellRange = $t.Cell(19,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "97-44="
$cellRange = $t.Cell(19,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "54+22="
$cellRange = $t.Cell(19,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "53+28="
$cellRange = $t.Cell(19,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "2+21="
$cellRange = $t.Cell(20,1).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "1+59="
$cellRange = $t.Cell(20,2).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "12+49="
$cellRange = $t.Cell(20,3).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "12+7="
$cellRange = $t.Cell(20,4).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "13+14="
$cellRange = $t.Cell(20,5).Range
[void]$cellRange.MoveEnd(1, -2)  # drop the trailing cell-end marks (\r\a)
$cellRange.Text = "20-4="
